$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look just like the existing header (AC1):
# bold, bordered, centered style. Copy the formatting from AC1 first, then
# set the text values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-43) gets the same team record: 77 wins, 85 losses, 0 ties.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 77
    $ws.Cells.Item($row, 31).Value = 85
    $ws.Cells.Item($row, 32).Value = 0
}
